# Update "想去人数" (F) and "最低票价" (G) values on the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1): row -> new F value
$sheet1F = @{
    2  = 126
    4  = 712
    8  = 2698
    9  = 1648
    10 = 1706
    12 = 278
    13 = 701
    15 = 142
    17 = 1107
    21 = 6090
    23 = 1251
    27 = 283
    29 = 53
    30 = 1078
    31 = 871
    33 = 78
    35 = 442
    36 = 1268
    37 = 153
    38 = 133
    41 = 158
}

foreach ($row in $sheet1F.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1F[$row]
}
# Also update the lowest ticket price (G) for row 31
$ws1.Cells.Item(31, 7).Value = 160

# Sheet "全部类型" (sheet4): row -> new F value
$sheet4F = @{
    2  = 126
    4  = 712
    11 = 2698
    12 = 1648
    13 = 1706
    15 = 278
    16 = 701
    19 = 142
    21 = 1107
    25 = 6090
    27 = 1251
    31 = 283
    33 = 53
    34 = 1078
    35 = 871
    37 = 78
    39 = 442
    40 = 1268
    41 = 153
    42 = 133
    45 = 158
}

foreach ($row in $sheet4F.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4F[$row]
}
# Also update the lowest ticket price (G) for row 35
$ws4.Cells.Item(35, 7).Value = 160
